$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '25.915.62'
Set-TextValue $ws.Range("E2") '  +0.01%  '

Set-TextValue $ws.Range("D3") '1.743.03'
Set-TextValue $ws.Range("E3") '  -0.60%  '

Set-TextValue $ws.Range("D4") '0.9998'
Set-TextValue $ws.Range("E4") '  -0.33%  '

Set-TextValue $ws.Range("D5") '230.63'
Set-TextValue $ws.Range("E5") '  -1.78%  '

Set-TextValue $ws.Range("D6") '0.9995'
Set-TextValue $ws.Range("E6") '  -0.24%  '

Set-TextValue $ws.Range("D7") '0.5256'
Set-TextValue $ws.Range("E7") '  +1.10%  '

Set-TextValue $ws.Range("D8") '0.2758'
Set-TextValue $ws.Range("E8") '  +1.43%  '

Set-TextValue $ws.Range("D9") '39.46'
Set-TextValue $ws.Range("E9") '  -2.45%  '

Set-TextValue $ws.Range("D10") '0.06151'
Set-TextValue $ws.Range("E10") '  +0.06%  '

Set-TextValue $ws.Range("D11") '1.739.18'
Set-TextValue $ws.Range("E11") '  -0.86%  '

Set-TextValue $ws.Range("E12") '  +0.86%  '

Set-TextValue $ws.Range("E13") '  -1.34%  '

Set-TextValue $ws.Range("D14") '0.6439'
Set-TextValue $ws.Range("E14") '  +1.72%  '

Set-TextValue $ws.Range("D15") '4.536'
Set-TextValue $ws.Range("E15") '  +0.54%  '

Set-TextValue $ws.Range("D16") '77.60'
Set-TextValue $ws.Range("E16") '  +0.12%  '

Set-TextValue $ws.Range("D17") '0.9998'
Set-TextValue $ws.Range("E17") '  -0.34%  '

Set-TextValue $ws.Range("D18") '0.9995'
Set-TextValue $ws.Range("E18") '  -0.23%  '

Set-TextValue $ws.Range("D19") '25.891.60'
Set-TextValue $ws.Range("E19") '  -0.13%  '

Set-TextValue $ws.Range("E20") '  -0.20%  '

Set-TextValue $ws.Range("D21") '0.000006691'
Set-TextValue $ws.Range("E21") '  +0.58%  '

Set-TextValue $ws.Range("D22") '1.963.39'
Set-TextValue $ws.Range("E22") '  -1.44%  '

Set-TextValue $ws.Range("D23") '4.304'
Set-TextValue $ws.Range("E23") '  +6.30%  '

Set-TextValue $ws.Range("D24") '8.784'
Set-TextValue $ws.Range("E24") '  +3.99%  '

Set-TextValue $ws.Range("D25") '5.166'
Set-TextValue $ws.Range("E25") '  +0.10%  '

Set-TextValue $ws.Range("D26") '140.26'
Set-TextValue $ws.Range("E26") '  +1.14%  '

Set-TextValue $ws.Range("E27") '  +0.73%  '

Set-TextValue $ws.Range("D28") '15.16'
Set-TextValue $ws.Range("E28") '  +0.59%  '

Set-TextValue $ws.Range("E29") '  -1.83%  '

Set-TextValue $ws.Range("D30") '102.71'
Set-TextValue $ws.Range("E30") '  -0.35%  '

Set-TextValue $ws.Range("D31") '0.08313'
Set-TextValue $ws.Range("E31") '  -0.38%  '

Set-TextValue $ws.Range("D32") '3.730'
Set-TextValue $ws.Range("E32") '  +2.17%  '

Set-TextValue $ws.Range("D33") '3.538'
Set-TextValue $ws.Range("E33") '  +3.56%  '

Set-TextValue $ws.Range("D34") '0.04532'
Set-TextValue $ws.Range("E34") '  +2.30%  '

Set-TextValue $ws.Range("D35") '2.616'
Set-TextValue $ws.Range("E35") '  -0.77%  '

Set-TextValue $ws.Range("D36") '0.9794'
Set-TextValue $ws.Range("E36") '  -1.13%  '

Set-TextValue $ws.Range("D37") '0.6228'
Set-TextValue $ws.Range("E37") '  +3.52%  '

Set-TextValue $ws.Range("D38") '2.684'
Set-TextValue $ws.Range("E38") '  -1.17%  '

Set-TextValue $ws.Range("E39") '  +0.62%  '

Set-TextValue $ws.Range("D40") '1.941'
Set-TextValue $ws.Range("E40") '  -0.07%  '

Set-TextValue $ws.Range("D41") '0.9992'
Set-TextValue $ws.Range("E41") '  -0.26%  '

Set-TextValue $ws.Range("D42") '100.26'
Set-TextValue $ws.Range("E42") '  -2.14%  '

Set-TextValue $ws.Range("D43") '0.3877'
Set-TextValue $ws.Range("E43") '  +0.85%  '

Set-TextValue $ws.Range("D44") '0.7319'
Set-TextValue $ws.Range("E44") '  -0.87%  '

Set-TextValue $ws.Range("D45") '5.009'
Set-TextValue $ws.Range("E45") '  +2.37%  '

Set-TextValue $ws.Range("D46") '0.05344'
Set-TextValue $ws.Range("E46") '  -2.98%  '

Set-TextValue $ws.Range("B47") 'Algorand'
Set-TextValue $ws.Range("C47") 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D47") '0.1129'
Set-TextValue $ws.Range("E47") '  +2.00%  '

Set-TextValue $ws.Range("B48") 'Aptos'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range("D48") '6.276'
Set-TextValue $ws.Range("E48") '  -0.05%  '

Set-TextValue $ws.Range("D49") '53.62'
Set-TextValue $ws.Range("E49") '  +2.72%  '

Set-TextValue $ws.Range("D50") '30.16'
Set-TextValue $ws.Range("E50") '  +0.58%  '

Set-TextValue $ws.Range("D51") '7.691'
Set-TextValue $ws.Range("E51") '  +3.68%  '
